# Add a new "Greece Market" test-data sheet, cloned from the "Croatia" sheet,
# matching the commit "Test data for Greece Market".

$wb = $excel.ActiveWorkbook

# Grab the last sheet (Croatia) - it's the template every market tab is based on.
$croatia = $wb.Worksheets.Item("Croatia")

# Activate it and select the whole sheet first, so its own tab records a
# "select all" state once it stops being the active tab (matches how Excel
# leaves the previously-active sheet after you duplicate+move off of it).
$croatia.Activate()
$croatia.Cells.Select()

# Duplicate "Croatia", inserting the copy right after it. Excel auto-assigns
# the next sheetId / rId and makes the new copy the active sheet.
$croatia.Copy($null, $croatia)

# The freshly created copy is now the last sheet in the workbook.
$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"

# Fill in the market-specific cells. Write B4 (ticket ref) before B2 (market
# name) so the two new shared-string entries land in the same order as the
# target workbook: "NGC-4119/T3165" then "Greece Market".
$greece.Range("B4").Value = "NGC-4119/T3165"
$greece.Range("B2").Value = "Greece Market"

# Leave the cursor where the author left it on the new sheet.
$greece.Range("B9").Select()
